$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Add the new row for "Saint Martin" edge case, right after the last used row (64)
$ws.Range("A65").Value = "Saint Martin"
$ws.Range("B65").Value = 1

# Match the resulting view: scrolled down a bit with F62 as the active selection
$ws.Range("F62").Select()

# Page setup now explicitly records portrait orientation
$ws.PageSetup.Orientation = 1
